# Insert a new weekly record for "Brócoli" (Macroferia Regional de Talca) as
# row 476 of the data sheet. Inserting the row shifts every following row
# down by one (old row 476 -> 477, ..., old row 511 -> 512), which matches
# the target diff: the sheet grows from A1:R511 to A1:R512 and the D/I/J/K/
# L/M/P values of rows 476-511 each become the values previously held by the
# row above them, while brand-new data appears in row 476.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 476..511 down to 477..512, duplicating formatting from row 476.
$ws.Rows(476).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(476, 1).Value  = 5
$ws.Cells.Item(476, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(476, 3).Value  = "Maule"
$ws.Cells.Item(476, 4).Value  = 45013
$ws.Cells.Item(476, 5).Value  = 7
$ws.Cells.Item(476, 6).Value  = 100112023
$ws.Cells.Item(476, 7).Value  = "Brócoli"
$ws.Cells.Item(476, 8).Value  = "Sin especificar"
$ws.Cells.Item(476, 9).Value  = "Primera"
$ws.Cells.Item(476, 10).Value = 5000
$ws.Cells.Item(476, 11).Value = 700
$ws.Cells.Item(476, 12).Value = 700
$ws.Cells.Item(476, 13).Value = 700
$ws.Cells.Item(476, 14).Value = "`$/unidad"
$ws.Cells.Item(476, 15).Value = "Región del Maule"
$ws.Cells.Item(476, 16).Value = 700
$ws.Cells.Item(476, 17).Value = 1
$ws.Cells.Item(476, 18).Value = "Hortaliza"
